# Update the "want to go" count (column F) for both the "展览" sheet
# and the "全部类型" sheet, which carry duplicate data.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 628
    3  = 577
    5  = 31
    6  = 124
    7  = 58
    8  = 58
    10 = 5055
    11 = 4722
    13 = 31
    16 = 178
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
